$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.546.98"
$ws.Range("E2").Value = "  +3.99%  "
$ws.Range("D3").Value = "2.451.25"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("E9").Value = "  +4.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.31%  "
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.66%  "
$ws.Range("D14").Value = "2.881.66"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").Value = "60.427.53"
$ws.Range("E15").Value = "  +3.87%  "
$ws.Range("E16").Value = "  +4.38%  "
$ws.Range("D17").Value = "2.450.71"
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.12%  "
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "0.0₃0800"
$ws.Range("E28").Value = "  +7.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("E35").Value = "  +5.18%  "
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("E40").Value = "  +10.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "317.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "144.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0966"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0527"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.575"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.404"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.72%  "
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.90%  "
